$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB price update
$ws.Range("D2").Value = "'229.42"

# Row 3 - OKB price update
$ws.Range("D3").Value = "'22.22"

# Row 4 - HuobiToken price update
$ws.Range("D4").Value = "'5.236"

# Row 5 - Cronos price update
$ws.Range("D5").Value = "'0.05558"

# Row 6 - GateToken price update
$ws.Range("D6").Value = "'3.377"

# Row 7 price update
$ws.Range("D7").Value = "'6.469"

# Row 8 price update
$ws.Range("D8").Value = "'1.058"

# Row 11 price update
$ws.Range("D11").Value = "'0.07320"

# Row 13 price update
$ws.Range("D13").Value = "'0.02948"

# Row 14 price update
$ws.Range("D14").Value = "'0.09270"

# Row 15 price update
$ws.Range("D15").Value = "'0.001656"

# Row 16 price update
$ws.Range("D16").Value = "'3.258"

# Row 17 price update
$ws.Range("D17").Value = "'0.04783"

# Row 18 - One price & worst-in-24h flag update
$ws.Range("D18").Value = "'0.0005888"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19 price update
$ws.Range("D19").Value = "'0.006221"

# Row 20 price update
$ws.Range("D20").Value = "'0.005230"

# Row 22 price update
$ws.Range("D22").Value = "'0.0001500"

# Row 23 price update
$ws.Range("D23").Value = "'3.905"

# Row 40 price update
$ws.Range("D40").Value = "'0.03995"

# Row 41 price update
$ws.Range("D41").Value = "'0.007132"

# Rows 42 & 43 - CEJI and BKEXToken swapped positions with new prices
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003501"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1036"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# Row 44 price update
$ws.Range("D44").Value = "'0.009944"

# Row 45 price update
$ws.Range("D45").Value = "'0.00005435"

# Row 47 price update
$ws.Range("D47").Value = "'0.7849"

# Row 48 - BOLO price update & worst-in-24h flag removed
$ws.Range("D48").Value = "'0.04297"
$ws.Range("E48").Value = "47BOLOBOLO"

# Row 49 price update
$ws.Range("D49").Value = "'0.00002100"
